# Add a new "2022-Q4" sheet (holdings detail) right after the "总计"
# sheet / before "2022-Q3", and add a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for 2022-Q4 and push the rest down.
# ---------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows("2:2").Insert()
$totals.Rows("2:2").ClearFormats()

# Re-use the existing "row marker" style (bold + border + centred) that
# already lives on A3 (the old A2) instead of re-building it by hand.
$totals.Range("A3").Copy($totals.Range("A2"))

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 10
$totals.Range("D2").Value = 0.1

# Renumber the index column for the rows that shifted down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3

# ---------------------------------------------------------------
# 2. Insert the brand-new "2022-Q4" worksheet before "2022-Q3".
# ---------------------------------------------------------------
$wb.Worksheets.Item("2022-Q3").Activate()
$q4 = $wb.Worksheets.Add()
$q4.Name = "2022-Q4"

# Match the page margins used by the sibling quarter sheets.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Copy the header row (same labels/style as every other quarter sheet)
# straight from "2022-Q3" so the styling matches exactly.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("B1:H1").Copy($q4.Range("B1:H1"))

# Broadcast the bold/border/centred "row marker" style used in column A
# down the whole data range.
$q3.Range("A2").Copy($q4.Range("A2:A11"))

# ---- data rows -------------------------------------------------
# Columns B-G are stored as text (matches the source workbook's
# convention of keeping these figures as strings); force text first so
# numeric-looking values ("0.89", "501305", ...) aren't coerced to
# numbers, then drop back to the default (unstyled) look.
$q4.Range("B2:G11").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "501305"
$q4.Range("C2").Value = "汇添富中证港股通高股息投资指数（LOF）A"
$q4.Range("D2").Value = "0.89"
$q4.Range("E2").Value = "91.24"
$q4.Range("F2").Value = "4.18"
$q4.Range("G2").Value = "0.0372"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "513530"
$q4.Range("C3").Value = "华泰柏瑞中证港股通高股息投资ETF（QDII）"
$q4.Range("D3").Value = "0.63"
$q4.Range("E3").Value = "96.34"
$q4.Range("F3").Value = "4.40"
$q4.Range("G3").Value = "0.0277"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "501306"
$q4.Range("C4").Value = "汇添富中证港股通高股息投资指数（LOF）C"
$q4.Range("D4").Value = "0.23"
$q4.Range("E4").Value = "91.24"
$q4.Range("F4").Value = "4.18"
$q4.Range("G4").Value = "0.0096"
$q4.Range("H4").Value = 3

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "004532"
$q4.Range("C5").Value = "民生加银中证港股通高股息精选指数A"
$q4.Range("D5").Value = "0.14"
$q4.Range("E5").Value = "92.86"
$q4.Range("F5").Value = "5.73"
$q4.Range("G5").Value = "0.0080"
$q4.Range("H5").Value = 4

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "004533"
$q4.Range("C6").Value = "民生加银中证港股通高股息精选指数C"
$q4.Range("D6").Value = "0.09"
$q4.Range("E6").Value = "92.86"
$q4.Range("F6").Value = "5.73"
$q4.Range("G6").Value = "0.0052"
$q4.Range("H6").Value = 4

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "006658"
$q4.Range("C7").Value = "财通中证香港红利等权投资指数A"
$q4.Range("D7").Value = "0.14"
$q4.Range("E7").Value = "89.84"
$q4.Range("F7").Value = "3.15"
$q4.Range("G7").Value = "0.0044"
$q4.Range("H7").Value = 5

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "501307"
$q4.Range("C8").Value = "银河中证沪港深高股息指数（LOF）A"
$q4.Range("D8").Value = "0.16"
$q4.Range("E8").Value = "93.15"
$q4.Range("F8").Value = "1.38"
$q4.Range("G8").Value = "0.0022"
$q4.Range("H8").Value = 7

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "006659"
$q4.Range("C9").Value = "财通中证香港红利等权投资指数C"
$q4.Range("D9").Value = "0.04"
$q4.Range("E9").Value = "89.84"
$q4.Range("F9").Value = "3.15"
$q4.Range("G9").Value = "0.0013"
$q4.Range("H9").Value = 5

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "005770"
$q4.Range("C10").Value = "信澳中证沪港深高股息精选指数"
$q4.Range("D10").Value = "0.13"
$q4.Range("E10").Value = "23.47"
$q4.Range("F10").Value = "0.63"
$q4.Range("G10").Value = "0.0008"
$q4.Range("H10").Value = 5

$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "501308"
$q4.Range("C11").Value = "银河中证沪港深高股息指数（LOF）C"
$q4.Range("D11").Value = "0.01"
$q4.Range("E11").Value = "93.15"
$q4.Range("F11").Value = "1.38"
$q4.Range("G11").Value = "0.0001"
$q4.Range("H11").Value = 7

# Drop the data cells back to the default (unstyled) look, matching the
# other quarter sheets, while keeping the values as text.
$q4.Range("B2:G11").Style = "Normal"

# ---------------------------------------------------------------
# 3. Restore the originally-selected tab ("2021-Q3" was the selected
#    sheet in the source workbook).
# ---------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
